$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row content (text in column A, frequency in column B, percentage in column C) ---
# Row 4: Amino acid metabolism
$ws.Range("A4").Value2 = "Amino acid metabolism"
$ws.Range("B4").Value2 = 30
$ws.Range("C4").Value2 = 30.612244897959183

# Row 5: Carbohydrate metabolism
$ws.Range("A5").Value2 = "Carbohydrate metabolism"
$ws.Range("B5").Value2 = 24
$ws.Range("C5").Value2 = 24.489795918367346

# Row 6: Nucleotide metabolism
$ws.Range("A6").Value2 = "Nucleotide metabolism"
$ws.Range("B6").Value2 = 19
$ws.Range("C6").Value2 = 19.387755102040817

# Row 7: Energy metabolism
$ws.Range("A7").Value2 = "Energy metabolism"
$ws.Range("B7").Value2 = 6
$ws.Range("C7").Value2 = 6.1224489795918364

# Row 8: Biosynthesis of other secondary metabolites (was Metabolism of cofactors and vitamins)
$ws.Range("A8").Value2 = "Biosynthesis of other secondary metabolites"
$ws.Range("B8").Value2 = 5
$ws.Range("C8").Value2 = 5.1020408163265305

# Row 9: Metabolism of other amino acids (was Biosynthesis of other secondary metabolites)
$ws.Range("A9").Value2 = "Metabolism of other amino acids"
$ws.Range("B9").Value2 = 4
$ws.Range("C9").Value2 = 4.0816326530612246

# Row 10: Metabolism of cofactors and vitamins (was Metabolism of other amino acids)
$ws.Range("A10").Value2 = "Metabolism of cofactors and vitamins"
$ws.Range("B10").Value2 = 4
$ws.Range("C10").Value2 = 4.0816326530612246

# Row 11: Glycan biosynthesis and metabolism
$ws.Range("A11").Value2 = "Glycan biosynthesis and metabolism"
$ws.Range("B11").Value2 = 3
$ws.Range("C11").Value2 = 3.0612244897959182

# Row 12: Metabolism of terpenoids and polyketides
$ws.Range("A12").Value2 = "Metabolism of terpenoids and polyketides"
$ws.Range("B12").Value2 = 2
$ws.Range("C12").Value2 = 2.0408163265306123

# Row 13: Lipid metabolism
$ws.Range("A13").Value2 = "Lipid metabolism"
$ws.Range("B13").Value2 = 1
$ws.Range("C13").Value2 = 1.0204081632653061

# --- Formatting: drop the centered "plain" style from A4:C13, then give C4:C13 a
#     one-decimal number format (no centering) ---
$ws.Range("A4:C13").ClearFormats()
$ws.Range("C4:C13").NumberFormat = "0.0"

# --- Selection moves from A2 to the newly formatted percentage column ---
$ws.Range("C4:C13").Select()
